$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '64.265.17'
$ws.Range('E2').Value = '  +4.41%  '
$ws.Range('D3').Value = '2.783.00'
$ws.Range('E3').Value = '  +4.65%  '
$ws.Range('D4').Value = "'0.998"
$ws.Range('E4').Value = '  -0.19%  '
$ws.Range('D5').Value = "'583.55"
$ws.Range('E5').Value = '  +0.19%  '
$ws.Range('D6').Value = "'161.41"
$ws.Range('E6').Value = '  +11.22%  '
$ws.Range('E7').Value = '  +3.50%  '
$ws.Range('E8').Value = '  -0.12%  '
$ws.Range('D9').Value = '2.806.18'
$ws.Range('E9').Value = '  +4.83%  '
$ws.Range('E10').Value = '  +3.33%  '
$ws.Range('E11').Value = '  +3.44%  '
$ws.Range('D12').Value = "'0.397"
$ws.Range('E12').Value = '  +3.77%  '
$ws.Range('E13').Value = '  +1.13%  '
$ws.Range('D14').Value = '3.274.19'
$ws.Range('E14').Value = '  +4.65%  '
$ws.Range('D15').Value = "'27.65"
$ws.Range('E15').Value = '  +5.85%  '
$ws.Range('D16').Value = '63.881.35'
$ws.Range('E16').Value = '  +4.00%  '
$ws.Range('E17').Value = '  +8.93%  '
$ws.Range('D18').Value = '2.790.62'
$ws.Range('E18').Value = '  +4.49%  '
$ws.Range('D19').Value = "'12.38"
$ws.Range('E19').Value = '  +5.72%  '
$ws.Range('E20').Value = '  +4.73%  '
$ws.Range('D21').Value = "'368.04"
$ws.Range('E21').Value = '  +3.06%  '
$ws.Range('D22').Value = "'7.09"
$ws.Range('E22').Value = '  +2.20%  '
$ws.Range('D23').Value = "'0.551"
$ws.Range('E23').Value = '  +4.94%  '
$ws.Range('E24').Value = '  +0.58%  '
$ws.Range('D25').Value = "'67.71"
$ws.Range('E25').Value = '  +4.58%  '
$ws.Range('D26').Value = "'0.175"
$ws.Range('E26').Value = '  +6.47%  '
$ws.Range('D27').Value = "'8.73"
$ws.Range('E27').Value = '  +2.57%  '
$ws.Range('D28').Value = '0.0₃0973'
$ws.Range('E28').Value = '  +17.36%  '
$ws.Range('D29').Value = "'0.997"
$ws.Range('E29').Value = '  +0.10%  '
$ws.Range('E30').Value = '  +1.36%  '
$ws.Range('D31').Value = "'7.36"
$ws.Range('E31').Value = '  +5.52%  '
$ws.Range('D32').Value = "'1.27"
$ws.Range('E32').Value = '  +11.63%  '
$ws.Range('D33').Value = "'173.86"
$ws.Range('E33').Value = '  +2.41%  '
$ws.Range('B34').Value = 'NEARProtocol'
$ws.Range('C34').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D34').Value = "'5.12"
$ws.Range('E34').Value = '  +8.98%  '
$ws.Range('B35').Value = 'EthereumClassic'
$ws.Range('C35').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D35').Value = "'20.95"
$ws.Range('E35').Value = '  +3.66%  '
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('D37').Value = "'1.50"
$ws.Range('E37').Value = '  +8.39%  '
$ws.Range('D38').Value = "'1.85"
$ws.Range('E39').Value = '  +3.64%  '
$ws.Range('B40').Value = 'Bittensor'
$ws.Range('C40').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D40').Value = "'342.98"
$ws.Range('E40').Value = '  -1.41%  '
$ws.Range('B41').Value = 'Filecoin'
$ws.Range('C41').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D41').Value = "'4.28"
$ws.Range('E41').Value = '  +3.08%  '
$ws.Range('D42').Value = "'6.26"
$ws.Range('E42').Value = '  +15.44%  '
$ws.Range('E43').Value = '  +3.08%  '
$ws.Range('D44').Value = "'22.65"
$ws.Range('E44').Value = '  +7.98%  '
$ws.Range('D45').Value = "'22.94"
$ws.Range('E45').Value = '  +7.69%  '
$ws.Range('D46').Value = "'0.0614"
$ws.Range('E46').Value = '  +5.62%  '
$ws.Range('E47').Value = '  +3.78%  '
$ws.Range('E48').Value = '  +2.98%  '
$ws.Range('D49').Value = "'138.68"
$ws.Range('E49').Value = '  +1.75%  '
$ws.Range('E50').Value = '  +2.55%  '
$ws.Range('D51').Value = '2.188.43'
$ws.Range('E51').Value = '  +3.96%  '
